# Update the "Marking" row (row 11) and "Total" row (row 12) on the
# "quiz" sheet to reflect the new concise_ms csv pattern.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Marking row: Right-answer weight 4 -> 5, Wrong-answer penalty -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Total row: Right total 44 -> 55, Wrong total 0 -> -0, Max fraction text updated
$ws.Range("B12").Value = 55
$ws.Range("C12").Value = -0
$ws.Range("E12").Value = "55.0/140"
